$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: "Rule For" ---

# Header J1: bold Arial header like the rest of row 1, but without the
# thick bottom border those header cells carry (matches the workbook's
# new cellXfs entry: same font as the header row, no border/fill).
$ws.Range("J1").Value = "Rule For"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J1").Borders.LineStyle = -4142   # xlLineStyleNone

# Data cells J2:J20: plain formatting, same as the existing column I data
# cells.
$ws.Range("I2").Copy()
$ws.Range("J2:J20").PasteSpecial(-4122)   # xlPasteFormats

$values = @(
    "Accounting", # J2
    "Reporting",  # J3
    "Accounting", # J4
    "Accounting", # J5
    "Accounting", # J6
    "Accounting", # J7
    "Accounting", # J8
    "Accounting", # J9
    "Accounting", # J10
    "Accounting", # J11
    "Accounting", # J12
    "Accounting", # J13
    "Accounting", # J14
    "Accounting", # J15
    "Accounting", # J16
    "Accounting", # J17
    "Accounting", # J18
    "Accounting", # J19
    "Accounting"  # J20
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $values[$i]
}

# Matches the selection left behind in the edited workbook.
$ws.Range("J4:J20").Select() | Out-Null
